# Applies numeric updates to the F ("想去人数") and G ("最低票价") columns
# across the four worksheets of the 上海-漫展信息 workbook, matching the
# commit "Update gh-pages to output generated at 456a3b4".
#
# NOTE: named parameters (e.g. "-Foo bar") are not reliably handled by this
# PowerShell runtime, so helper functions below are called positionally.

$wb = $excel.ActiveWorkbook

function SetCells {
    param($SheetName, $Updates)
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($cellRef in $Updates.Keys) {
        $ws.Range($cellRef).Value = $Updates[$cellRef]
    }
}

SetCells "展览" @{
    "F4"  = 8441
    "F6"  = 154
    "F7"  = 2368
    "F9"  = 132
    "F13" = 1032
    "G13" = 75
    "F14" = 1585
    "F15" = 2183
    "F17" = 235
    "F18" = 299
    "F19" = 2041
    "F21" = 816
    "F22" = 800
    "F24" = 814
    "F25" = 1381
    "F26" = 579
    "F27" = 1209
    "F28" = 28
    "F29" = 271
    "F30" = 43
    "F33" = 2551
}

SetCells "演出" @{
    "F2"  = 137
    "F6"  = 28
    "F24" = 3
    "F27" = 27
    "F29" = 108
    "F32" = 7
    "F38" = 10
    "F39" = 277
    "F43" = 4
    "F48" = 38
}

SetCells "本地生活" @{
    "F4"  = 2460
    "F7"  = 706
    "F8"  = 2459
    "F9"  = 9514
    "F15" = 335
    "F16" = 2682
    "F17" = 323
    "F18" = 160
    "F19" = 615
}

SetCells "全部类型" @{
    "F4"  = 2460
    "F5"  = 706
    "F8"  = 154
    "F9"  = 335
    "F10" = 323
    "F11" = 132
    "F13" = 160
    "F14" = 1032
    "G14" = 75
    "F15" = 1585
    "F16" = 137
    "F17" = 615
    "F18" = 615
    "F20" = 28
    "F24" = 299
    "F25" = 2041
    "F29" = 814
    "F31" = 1381
    "F35" = 579
    "F38" = 3
    "F39" = 1210
    "F40" = 271
    "F41" = 27
    "F46" = 277
    "F47" = 2551
    "F49" = 4
}
